$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.402001857757568
$ws.Range("E2").Value = 1111.489525668081
$ws.Range("F2").Value = 0.04382238405149886
$ws.Range("G2").Value = 0.03605672573600956
$ws.Range("H2").Value = 0.03238454828618866
$ws.Range("I2").Value = 0.02984968898429384
$ws.Range("J2").Value = 0.02794665090691812
$ws.Range("K2").Value = 0.02685866724019624
$ws.Range("L2").Value = 0.02637038412424178
$ws.Range("M2").Value = 0.02559222481632079
$ws.Range("N2").Value = 0.02496187254004978
$ws.Range("O2").Value = 0.02424352554610183
$ws.Range("P2").Value = 0.02366113582582776
$ws.Range("Q2").Value = 0.02352623618664306
$ws.Range("R2").Value = 0.0229729996708121
$ws.Range("S2").Value = 0.02229779156302234
$ws.Range("T2").Value = 0.02229779156302234
$ws.Range("U2").Value = 0.02198308584230929
$ws.Range("V2").Value = 0.02191178504773458
$ws.Range("W2").Value = 0.0218204920310211
$ws.Range("X2").Value = 0.02172574359155547
$ws.Range("Y2").Value = 0.02166646248865655

$ws.Range("C3").Value = 1.236000299453735
$ws.Range("E3").Value = 1082.860336986898
$ws.Range("F3").Value = 0.04439805636963341
$ws.Range("G3").Value = 0.03764234423516215
$ws.Range("H3").Value = 0.03176317813667204
$ws.Range("I3").Value = 0.02952554752304999
$ws.Range("J3").Value = 0.0279123239030453
$ws.Range("K3").Value = 0.02656049370387836
$ws.Range("L3").Value = 0.02656049370387836
$ws.Range("M3").Value = 0.02555220208207397
$ws.Range("N3").Value = 0.02519183101446608
$ws.Range("O3").Value = 0.02405296854434908
$ws.Range("P3").Value = 0.02344602903885229
$ws.Range("Q3").Value = 0.0228845022980148
$ws.Range("R3").Value = 0.02208599530848862
$ws.Range("S3").Value = 0.02205893878578527
$ws.Range("T3").Value = 0.02179551034581526
$ws.Range("U3").Value = 0.02153673456125409
$ws.Range("V3").Value = 0.02134914433113034
$ws.Range("W3").Value = 0.02115645459621594
$ws.Range("X3").Value = 0.02115645459621594
$ws.Range("Y3").Value = 0.02110838863522218

$ws.Range("C4").Value = 1.078995227813721
$ws.Range("E4").Value = 1092.982650034321
$ws.Range("F4").Value = 0.04324704218586918
$ws.Range("G4").Value = 0.03560538235314089
$ws.Range("H4").Value = 0.03193136677862912
$ws.Range("I4").Value = 0.03068243813614821
$ws.Range("J4").Value = 0.02896938807578996
$ws.Range("K4").Value = 0.02744447007205675
$ws.Range("L4").Value = 0.02669171560544902
$ws.Range("M4").Value = 0.02526982430849956
$ws.Range("N4").Value = 0.02467936451253348
$ws.Range("O4").Value = 0.02407866592412512
$ws.Range("P4").Value = 0.02331085435339228
$ws.Range("Q4").Value = 0.02268952003632042
$ws.Range("R4").Value = 0.02225602004546522
$ws.Range("S4").Value = 0.02209841139334741
$ws.Range("T4").Value = 0.02197846537064103
$ws.Range("U4").Value = 0.02179945461747924
$ws.Range("V4").Value = 0.02142021076992367
$ws.Range("W4").Value = 0.02142021076992367
$ws.Range("X4").Value = 0.02130570467903159
$ws.Range("Y4").Value = 0.02130570467903159

$ws.Range("C5").Value = 1.214994192123413
$ws.Range("E5").Value = 1098.432134141374
$ws.Range("F5").Value = 0.0430273248552243
$ws.Range("G5").Value = 0.03637295598952074
$ws.Range("H5").Value = 0.0324591586329178
$ws.Range("I5").Value = 0.02918095392931441
$ws.Range("J5").Value = 0.02830039816518203
$ws.Range("K5").Value = 0.02690296354921926
$ws.Range("L5").Value = 0.02573810522365879
$ws.Range("M5").Value = 0.02507692882594763
$ws.Range("N5").Value = 0.02499462474191791
$ws.Range("O5").Value = 0.02326354696619991
$ws.Range("P5").Value = 0.02326354696619991
$ws.Range("Q5").Value = 0.02301018492677821
$ws.Range("R5").Value = 0.02261881099484565
$ws.Range("S5").Value = 0.02237915874756227
$ws.Range("T5").Value = 0.0220023462817038
$ws.Range("U5").Value = 0.02181281731084303
$ws.Range("V5").Value = 0.02181281731084303
$ws.Range("W5").Value = 0.02159162198020709
$ws.Range("X5").Value = 0.02157069656937488
$ws.Range("Y5").Value = 0.021411932439403

$ws.Range("C6").Value = 1.074001550674438
$ws.Range("E6").Value = 1109.690256389744
$ws.Range("F6").Value = 0.04322574068151943
$ws.Range("G6").Value = 0.03602517710291294
$ws.Range("H6").Value = 0.03366921202197963
$ws.Range("I6").Value = 0.03029766197714675
$ws.Range("J6").Value = 0.02805677936716721
$ws.Range("K6").Value = 0.02614089965699919
$ws.Range("L6").Value = 0.02592307871928971
$ws.Range("M6").Value = 0.02408444495740039
$ws.Range("N6").Value = 0.02375540266394418
$ws.Range("O6").Value = 0.02375540266394418
$ws.Range("P6").Value = 0.02344053090514753
$ws.Range("Q6").Value = 0.02342130949597616
$ws.Range("R6").Value = 0.02322807208097031
$ws.Range("S6").Value = 0.02253845672827561
$ws.Range("T6").Value = 0.02253845672827561
$ws.Range("U6").Value = 0.02222802160163564
$ws.Range("V6").Value = 0.02193913932402497
$ws.Range("W6").Value = 0.02180207921234046
$ws.Range("X6").Value = 0.0217133420460938
$ws.Range("Y6").Value = 0.02163138901344529

$ws.Range("C7").Value = 1.234985828399658
$ws.Range("E7").Value = 1077.53029752072
$ws.Range("F7").Value = 0.04441323764315241
$ws.Range("G7").Value = 0.03718071724471769
$ws.Range("H7").Value = 0.0329524146774493
$ws.Range("I7").Value = 0.03161868548943408
$ws.Range("J7").Value = 0.02797589206553492
$ws.Range("K7").Value = 0.02758838242935261
$ws.Range("L7").Value = 0.0261139364447245
$ws.Range("M7").Value = 0.02554467800895161
$ws.Range("N7").Value = 0.0251457368864357
$ws.Range("O7").Value = 0.0239292109474071
$ws.Range("P7").Value = 0.02328091686384322
$ws.Range("Q7").Value = 0.02273191527302573
$ws.Range("R7").Value = 0.02217539282607406
$ws.Range("S7").Value = 0.02193018481486422
$ws.Range("T7").Value = 0.02142155465755702
$ws.Range("U7").Value = 0.02142155465755702
$ws.Range("V7").Value = 0.02130395643118798
$ws.Range("W7").Value = 0.02114923592997328
$ws.Range("X7").Value = 0.02109244753830948
$ws.Range("Y7").Value = 0.02100448923042338

$ws.Range("C8").Value = 1.129998683929443
$ws.Range("E8").Value = 1080.060117035959
$ws.Range("F8").Value = 0.04374821058621453
$ws.Range("G8").Value = 0.03578353419861968
$ws.Range("H8").Value = 0.03311995808889262
$ws.Range("I8").Value = 0.02975908129897101
$ws.Range("J8").Value = 0.02870288693477724
$ws.Range("K8").Value = 0.02660013566609448
$ws.Range("L8").Value = 0.02546180130593555
$ws.Range("M8").Value = 0.02456785495300525
$ws.Range("N8").Value = 0.02437522386763857
$ws.Range("O8").Value = 0.02364254696633713
$ws.Range("P8").Value = 0.02329047287497531
$ws.Range("Q8").Value = 0.02255923391821012
$ws.Range("R8").Value = 0.02208298668887341
$ws.Range("S8").Value = 0.02195256541066474
$ws.Range("T8").Value = 0.02186513080063413
$ws.Range("U8").Value = 0.02156521346846306
$ws.Range("V8").Value = 0.02149899777680111
$ws.Range("W8").Value = 0.02133298732568251
$ws.Range("X8").Value = 0.02114169244952361
$ws.Range("Y8").Value = 0.02105380345099334

$ws.Range("C9").Value = 1.120998859405518
$ws.Range("E9").Value = 1117.965144351105
$ws.Range("F9").Value = 0.04343125467729296
$ws.Range("G9").Value = 0.03641943870275671
$ws.Range("H9").Value = 0.03254155954424213
$ws.Range("I9").Value = 0.02991022402329768
$ws.Range("J9").Value = 0.0280023202606109
$ws.Range("K9").Value = 0.02738840712066392
$ws.Range("L9").Value = 0.02580314601145084
$ws.Range("M9").Value = 0.02532194714767779
$ws.Range("N9").Value = 0.0247402815112107
$ws.Range("O9").Value = 0.02391251089834247
$ws.Range("P9").Value = 0.02376958059690317
$ws.Range("Q9").Value = 0.02347955588099854
$ws.Range("R9").Value = 0.02330064599207457
$ws.Range("S9").Value = 0.0226390589394497
$ws.Range("T9").Value = 0.02244338473472619
$ws.Range("U9").Value = 0.02244338473472619
$ws.Range("V9").Value = 0.02210952637079219
$ws.Range("W9").Value = 0.02205329238684126
$ws.Range("X9").Value = 0.02181332420809654
$ws.Range("Y9").Value = 0.02179269287234122

$ws.Range("C10").Value = 1.188011884689331
$ws.Range("E10").Value = 1117.0367702164
$ws.Range("F10").Value = 0.04288414632950576
$ws.Range("G10").Value = 0.03653695853235905
$ws.Range("H10").Value = 0.03258139865763486
$ws.Range("I10").Value = 0.03153219665857632
$ws.Range("J10").Value = 0.03045154639991922
$ws.Range("K10").Value = 0.02897521860848758
$ws.Range("L10").Value = 0.02634211980400364
$ws.Range("M10").Value = 0.0263360261240581
$ws.Range("N10").Value = 0.0254773467478316
$ws.Range("O10").Value = 0.02415855973949551
$ws.Range("P10").Value = 0.02396817250775325
$ws.Range("Q10").Value = 0.02368186948899537
$ws.Range("R10").Value = 0.02299909999260738
$ws.Range("S10").Value = 0.02266921818565666
$ws.Range("T10").Value = 0.02241552012280473
$ws.Range("U10").Value = 0.02215916853157435
$ws.Range("V10").Value = 0.02215916853157435
$ws.Range("W10").Value = 0.02215593922634545
$ws.Range("X10").Value = 0.02195319362395484
$ws.Range("Y10").Value = 0.02177459591065107

$ws.Range("C11").Value = 1.062998533248901
$ws.Range("E11").Value = 1093.99302506202
$ws.Range("F11").Value = 0.04380072938683001
$ws.Range("G11").Value = 0.03747125524319624
$ws.Range("H11").Value = 0.03439291722931721
$ws.Range("I11").Value = 0.03087919643449778
$ws.Range("J11").Value = 0.02991392602388328
$ws.Range("K11").Value = 0.0270714242002147
$ws.Range("L11").Value = 0.02643852547618391
$ws.Range("M11").Value = 0.02476907238523677
$ws.Range("N11").Value = 0.02390815163129871
$ws.Range("O11").Value = 0.02386898097615687
$ws.Range("P11").Value = 0.0238288593150442
$ws.Range("Q11").Value = 0.02334602946096209
$ws.Range("R11").Value = 0.02239942446104278
$ws.Range("S11").Value = 0.02239942446104278
$ws.Range("T11").Value = 0.02203512979043112
$ws.Range("U11").Value = 0.02194339384576251
$ws.Range("V11").Value = 0.02177449649752153
$ws.Range("W11").Value = 0.0215414733934439
$ws.Range("X11").Value = 0.0214357854117811
$ws.Range("Y11").Value = 0.02132540009867485
